$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "hello"
$ws.Range("B1").Value = "code"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "mn"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "mnn"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "mnnn"

$ws.Range("B5").Select() | Out-Null
